$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.02%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.72"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.21%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.711"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.24%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06204"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.83%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.720"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.70%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8500"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.47%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9129"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.33%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.15%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04982"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.11%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07106"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.03%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03118"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.17%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09055"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001544"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.76%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006158"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.30%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005954"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.82%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.448"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.02%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.173"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.05%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.23%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1310"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.00%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.096"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04237"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.20%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001183"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.37%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004059"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.77%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.10%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.12%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03937"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.59%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1114"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.05%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004123"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.78%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.31%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-18.49%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005166"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.10%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "83.49%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.10%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.10%"
